$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 2348.614
$ws.Range("I15").Value = 2348.614
$ws.Range("K15").Value = 7045.842000000001
$ws.Range("M15").Value = -6876.842000000001
# row 51
$ws.Range("H51").Value = 4032.8
$ws.Range("J51").Value = 4408.727
$ws.Range("L51").Value = 4408.727
$ws.Range("N51").Value = -5376.727
# row 70
$ws.Range("H70").Value = 27766.143
$ws.Range("J70").Value = 28227.166
$ws.Range("L70").Value = 84681.49800000001
$ws.Range("N70").Value = -85221.49800000001
# row 73
$ws.Range("H73").Value = 27766.143
$ws.Range("J73").Value = 28227.166
$ws.Range("L73").Value = 84681.49800000001
$ws.Range("N73").Value = -86553.49800000001
# row 76
$ws.Range("H76").Value = 4014.6667
$ws.Range("J76").Value = 5500
$ws.Range("L76").Value = 5500
$ws.Range("N76").Value = -6130
# row 79
$ws.Range("H79").Value = 4014.6667
$ws.Range("J79").Value = 5500
$ws.Range("L79").Value = 5500
$ws.Range("N79").Value = -7684
# row 86
$ws.Range("H86").Value = 1009.53845
$ws.Range("I86").Value = 1218.1428
$ws.Range("J86").Value = 766.1667
$ws.Range("K86").Value = 1218.1428
$ws.Range("L86").Value = 766.1667
$ws.Range("M86").Value = -95.14280000000008
$ws.Range("N86").Value = -3012.1667
# row 88
$ws.Range("H88").Value = 6714.7144
$ws.Range("I88").Value = 1001.5
$ws.Range("K88").Value = 1001.5
$ws.Range("M88").Value = -595.5
# row 89
$ws.Range("H89").Value = 1009.53845
$ws.Range("I89").Value = 1218.1428
$ws.Range("J89").Value = 766.1667
$ws.Range("K89").Value = 6090.714
$ws.Range("L89").Value = 3830.8335
$ws.Range("M89").Value = -474.7139999999999
$ws.Range("N89").Value = -15062.8335
# row 91
$ws.Range("H91").Value = 6714.7144
$ws.Range("I91").Value = 1001.5
$ws.Range("K91").Value = 1001.5
$ws.Range("M91").Value = 402.5
# row 100
$ws.Range("H100").Value = 1971.5333
$ws.Range("I100").Value = 1300.6364
$ws.Range("J100").Value = 3816.5
$ws.Range("K100").Value = 1300.6364
$ws.Range("L100").Value = 3816.5
$ws.Range("M100").Value = -759.6364000000001
$ws.Range("N100").Value = -4898.5
# row 107
$ws.Range("H107").Value = 1048
$ws.Range("I107").Value = 660.125
$ws.Range("K107").Value = 660.125
$ws.Range("M107").Value = 1259.875
# row 137
$ws.Range("H137").Value = 26680.975
$ws.Range("I137").Value = 778.7646999999999
$ws.Range("J137").Value = 45826.086
$ws.Range("K137").Value = 2336.2941
$ws.Range("L137").Value = 137478.258
$ws.Range("M137").Value = 213.7058999999999
$ws.Range("N137").Value = -142578.258

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 3324087.5
$ws.Range("I2").Value = 4652722.5
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 4652722.5
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = -4652609.5
$ws.Range("N2").Value = -2726
# row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
# row 32
$ws.Range("H32").Value = 3020.982
$ws.Range("I32").Value = 2199.186
$ws.Range("K32").Value = 2199.186
$ws.Range("M32").Value = -1912.186
# row 37
$ws.Range("H37").Value = 16780
$ws.Range("J37").Value = 17725
$ws.Range("L37").Value = 17725
$ws.Range("N37").Value = -18271
# row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# row 61
$ws.Range("H61").Value = 4441.4
$ws.Range("I61").Value = 1700
$ws.Range("J61").Value = 4746
$ws.Range("K61").Value = 1700
$ws.Range("L61").Value = 4746
$ws.Range("M61").Value = -1488
$ws.Range("N61").Value = -5170
# row 74
$ws.Range("H74").Value = 1274.8
$ws.Range("I74").Value = 992.5714
$ws.Range("J74").Value = 1933.3334
$ws.Range("K74").Value = 992.5714
$ws.Range("L74").Value = 1933.3334
$ws.Range("M74").Value = -118.5714
$ws.Range("N74").Value = -3681.3334
# row 77
$ws.Range("H77").Value = 1274.8
$ws.Range("I77").Value = 992.5714
$ws.Range("J77").Value = 1933.3334
$ws.Range("K77").Value = 4962.857
$ws.Range("L77").Value = 9666.666999999999
$ws.Range("M77").Value = -594.857
$ws.Range("N77").Value = -18402.667
# row 97
$ws.Range("H97").Value = 2336.3333
$ws.Range("I97").Value = 2336.3333
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2336.3333
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1840.3333
$ws.Range("N97").ClearContents()
# row 116
$ws.Range("H116").Value = 3324087.5
$ws.Range("I116").Value = 4652722.5
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 4652722.5
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = -4650428.5
$ws.Range("N116").Value = -7088
# row 136
$ws.Range("H136").Value = 4441.4
$ws.Range("I136").Value = 1700
$ws.Range("J136").Value = 4746
$ws.Range("K136").Value = 5100
$ws.Range("L136").Value = 14238
$ws.Range("M136").Value = -2550
$ws.Range("N136").Value = -19338

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 3324087.5
$ws.Range("I3").Value = 4652722.5
$ws.Range("J3").Value = 2500
$ws.Range("K3").Value = 4652722.5
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = -4652608.5
$ws.Range("N3").Value = -2728
# row 20
$ws.Range("H20").Value = 1962.2222
$ws.Range("I20").Value = 1857.4736
$ws.Range("K20").Value = 1857.4736
$ws.Range("M20").Value = -1610.4736

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 4808272
$ws.Range("I22").Value = 451.16666
$ws.Range("J22").Value = 8929261
$ws.Range("K22").Value = 451.16666
$ws.Range("L22").Value = 8929261
$ws.Range("M22").Value = -101.16666
$ws.Range("N22").Value = -8929961
# row 50
$ws.Range("H50").Value = 17966.666
$ws.Range("J50").Value = 17966.666
$ws.Range("L50").Value = 17966.666
$ws.Range("N50").Value = -19216.666
# row 58
$ws.Range("I58").Value = 6212755
$ws.Range("J58").Value = 2147.8
$ws.Range("K58").Value = 6212755
$ws.Range("L58").Value = 2147.8
$ws.Range("M58").Value = -6212552
$ws.Range("N58").Value = -2553.8
# row 134
$ws.Range("H134").Value = 1595.4
$ws.Range("I134").Value = 1063
$ws.Range("J134").Value = 3725
$ws.Range("K134").Value = 3189
$ws.Range("L134").Value = 11175
$ws.Range("M134").Value = -654
$ws.Range("N134").Value = -16245
# row 136
$ws.Range("I136").Value = 6212755
$ws.Range("J136").Value = 2147.8
$ws.Range("K136").Value = 18638265
$ws.Range("L136").Value = 6443.400000000001
$ws.Range("M136").Value = -18635715
$ws.Range("N136").Value = -11543.4

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 61
$ws.Range("H61").Value = 199.33333
$ws.Range("J61").Value = 198
$ws.Range("L61").Value = 594
$ws.Range("N61").Value = -1024
# row 113
$ws.Range("H113").Value = 111810.664
$ws.Range("I113").Value = 1000003
$ws.Range("J113").Value = 786.625
$ws.Range("K113").Value = 3000009
$ws.Range("L113").Value = 2359.875
$ws.Range("M113").Value = -2997839
$ws.Range("N113").Value = -6699.875
# row 131
$ws.Range("H131").Value = 12840170
$ws.Range("J131").Value = 21850.514
$ws.Range("L131").Value = 65551.542
$ws.Range("N131").Value = -75631.542
# row 140
$ws.Range("H140").Value = 2961.0278
$ws.Range("I140").Value = 609.55554
$ws.Range("K140").Value = 1828.66662
$ws.Range("M140").Value = 3351.33338

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 240.3
$ws.Range("I2").Value = 255.85715
$ws.Range("K2").Value = 255.85715
$ws.Range("M2").Value = -142.85715
# row 132
$ws.Range("H132").Value = 1541484.4
$ws.Range("I132").Value = 2264188
$ws.Range("J132").Value = 5738.875
$ws.Range("K132").Value = 6792564
$ws.Range("L132").Value = 17216.625
$ws.Range("M132").Value = -6790034
$ws.Range("N132").Value = -22276.625
# row 141
$ws.Range("H141").Value = 39500
$ws.Range("J141").Value = 39500
$ws.Range("L141").Value = 39500
$ws.Range("N141").Value = -49860

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 3828.2778
$ws.Range("J7").Value = 4980.5
$ws.Range("L7").Value = 4980.5
$ws.Range("N7").Value = -5204.5
# row 46
$ws.Range("H46").Value = 1190.8334
$ws.Range("I46").Value = 672
$ws.Range("J46").Value = 1363.7778
$ws.Range("K46").Value = 672
$ws.Range("L46").Value = 1363.7778
$ws.Range("M46").Value = -484
$ws.Range("N46").Value = -1739.7778
# row 68
$ws.Range("H68").Value = 3781.4285
$ws.Range("I68").Value = 3578.3333
$ws.Range("K68").Value = 3578.3333
$ws.Range("M68").Value = -2829.3333
# row 71
$ws.Range("H71").Value = 3781.4285
$ws.Range("I71").Value = 3578.3333
$ws.Range("K71").Value = 17891.6665
$ws.Range("M71").Value = -14147.6665
# row 126
$ws.Range("H126").Value = 3828.2778
$ws.Range("J126").Value = 4980.5
$ws.Range("L126").Value = 14941.5
$ws.Range("N126").Value = -19881.5
# row 132
$ws.Range("H132").Value = 3123.1724
$ws.Range("I132").Value = 1439
$ws.Range("K132").Value = 4317
$ws.Range("M132").Value = -1787

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 51
$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4490
# row 104
$ws.Range("H104").Value = 15000
$ws.Range("J104").Value = 15000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -21988
# row 132
$ws.Range("H132").Value = 4999
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
# row 135
$ws.Range("H135").Value = 90282.60000000001
$ws.Range("J135").Value = 90282.60000000001
$ws.Range("L135").Value = 90282.60000000001
$ws.Range("N135").Value = -100422.6
